$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "KZtcy695"
$ws.Range("B2").Value = 23081625
$ws.Range("C2").Value = "bdcgwbt47"
$ws.Range("D2").Value = "hJk7!3#R"
$ws.Range("F2").Value = "jjcBWIZR"
$ws.Range("G2").Value = "hukc"
